# Pedido 69134d11b9c1d30b15fabdc3
# The sheet currently has a single data row (row 2) for the "Test Ringover
# (NO TOCAR)" order. This edit duplicates that row into a new row 3 (an
# identical second line item for the same order) and tidies up row 2 by
# dropping the placeholder/blank columns (Paneles, Unidades Paneles,
# Optimizador, Unidades Optimizador, Cargador VE, Unidades Cargador VE,
# Pajareras, Unidades Pajareras, LEG) that never had any real content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Duplicate row 2 (the whole data row, as it exists today) into row 3.
#    Using a values-only paste keeps the original text representation of
#    values such as "1" (stored as text, not numbers) without carrying any
#    extra formatting along with the copy.
$ws.Range("A2:R2").Copy()
$ws.Range("A3").PasteSpecial(-4163, -4142, $false, $false)
$excel.CutCopyMode = $false

# 2) Remove the empty placeholder cells from row 2 (Paneles/Optimizador/
#    Cargador VE/Pajareras/LEG columns), which were never used for this
#    order.
$ws.Range("E2:H2").ClearContents()
$ws.Range("M2:P2").ClearContents()
$ws.Range("R2").ClearContents()
